# "API register, verify, login"
# Update planning figures for the "Express, authentification" task (row 45)
# and related rows; dependent SUM/ratio formulas recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 38 - "Mise en place d'Express": réalisées hours 4 -> 6 (F38 = E38/D38 recalculates)
$ws.Range("E38").Value = 6

# Row 45 - "Express, authentification": planifiées 8 -> 16, réalisées 9 -> 13
$ws.Range("D45").Value = 16
$ws.Range("E45").Value = 13

# Row 49 - "À faire" / Réalisation: planifiées hours 200 -> 185
$ws.Range("D49").Value = 185

# Update the last active selection to reflect where the author finished editing
$ws.Range("K48").Select()
